# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-looking string that must remain TEXT
# (matches the source data which stores all Price-column entries as text).
# Setting NumberFormat to "@" (Text) before assigning the value prevents Excel
# from silently converting these into numeric cells.
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D15",
    "D17",
    "D20",
    "D21",
    "D22",
    "D23",
    "D27",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D41",
    "D44",
    "D45",
    "D46",
    "D47",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "41.353.24"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "2.440.22"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "317.43"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").Value = "89.43"
$ws.Range("E6").Value = "  -4.04%  "
$ws.Range("E7").Value = "  -2.29%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  -4.22%  "
$ws.Range("D10").Value = "32.14"
$ws.Range("E10").Value = "  -2.85%  "
$ws.Range("D11").Value = "0.0830"
$ws.Range("E11").Value = "  -6.96%  "
$ws.Range("D12").Value = "0.109"
$ws.Range("E12").Value = "  -2.65%  "
$ws.Range("D13").Value = "2.812.53"
$ws.Range("E13").Value = "  -1.55%  "
$ws.Range("D14").Value = "6.72"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("D15").Value = "15.51"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "2.461.99"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").Value = "0.774"
$ws.Range("E17").Value = "  -2.39%  "
$ws.Range("D18").Value = "41.265.08"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "0.0₃0922"
$ws.Range("E19").Value = "  -4.20%  "
$ws.Range("D20").Value = "6.24"
$ws.Range("E20").Value = "  -3.93%  "
$ws.Range("D21").Value = "71.78"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "11.07"
$ws.Range("E22").Value = "  -3.85%  "
$ws.Range("D23").Value = "235.41"
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("E24").Value = "  -2.27%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("D27").Value = "24.06"
$ws.Range("E27").Value = "  -3.55%  "
$ws.Range("E28").Value = "  -3.39%  "
$ws.Range("D29").Value = "9.56"
$ws.Range("E29").Value = "  -3.50%  "
$ws.Range("D30").Value = "34.69"
$ws.Range("E30").Value = "  -4.95%  "
$ws.Range("D31").Value = "157.61"
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "5.27"
$ws.Range("E32").Value = "  -4.92%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "2.53"
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("D35").Value = "0.0745"
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("D36").Value = "2.91"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "16.58"
$ws.Range("E37").Value = "  -5.47%  "
$ws.Range("D38").Value = "0.115"
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("E39").Value = "  -3.51%  "
$ws.Range("E40").Value = "  -3.83%  "
$ws.Range("D41").Value = "3.89"
$ws.Range("E41").Value = "  -2.70%  "
$ws.Range("E42").Value = "  -7.21%  "
$ws.Range("D43").Value = "1.985.85"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "0.0275"
$ws.Range("E44").Value = "  -3.84%  "
$ws.Range("D45").Value = "18.23"
$ws.Range("E45").Value = "  -5.96%  "
$ws.Range("D46").Value = "2.87"
$ws.Range("E46").Value = "  -4.93%  "
$ws.Range("D47").Value = "9.53"
$ws.Range("E47").Value = "  +3.91%  "
$ws.Range("D48").Value = "2.671.36"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D49").Value = "95.17"
$ws.Range("E49").Value = "  -2.53%  "
$ws.Range("D50").Value = "73.18"
$ws.Range("E50").Value = "  -1.27%  "
$ws.Range("D51").Value = "52.06"
$ws.Range("E51").Value = "  -1.48%  "
